# Weekly data refresh: a new week's price record is inserted at row 25,
# pushing all subsequent rows (old 25..78) down by one (new 26..79).
# The sheet's dimension grows from A1:R78 to A1:R79.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 25 (shifts rows 25:78 down to 26:79,
# carrying formatting from the row above - matches native Excel behaviour).
$ws.Rows.Item(25).Insert()

# Populate the newly inserted row 25 with the new weekly record.
$ws.Cells.Item(25, 1).Value  = 11
$ws.Cells.Item(25, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(25, 3).Value  = "Bíobío"
$ws.Cells.Item(25, 4).Value  = 45028
$ws.Cells.Item(25, 5).Value  = 8
$ws.Cells.Item(25, 6).Value  = 100112031
$ws.Cells.Item(25, 7).Value  = "Poroto verde"
$ws.Cells.Item(25, 8).Value  = "Magnum"
$ws.Cells.Item(25, 9).Value  = "Primera"
$ws.Cells.Item(25, 10).Value = 110
$ws.Cells.Item(25, 11).Value = 14500
$ws.Cells.Item(25, 12).Value = 15000
$ws.Cells.Item(25, 13).Value = 14727
$ws.Cells.Item(25, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(25, 15).Value = "Región del Maule"
$ws.Cells.Item(25, 16).Value = 589
$ws.Cells.Item(25, 17).Value = 25
$ws.Cells.Item(25, 18).Value = "Hortaliza"
